$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the trailing placeholder rows (1048571-1048576) that only carried a
# row height with no cell content -- the final sheet doesn't keep them.
$ws.Range("A1048571:A1048576").EntireRow.Delete()

# Insert 5 new rows below the existing data (rows 7-11), copying the format
# (style 4 for A:C, style 1 for D) from the last data row (row 2) so the
# inserted rows inherit the same look as the "eng" rows.
$ws.Range("A2:D2").Copy()
$ws.Range("A7:A11").EntireRow.Insert(-4121)

# Fill in the new French ("fra") authentication-method rows.
$ws.Range("A7").Value = "fra"
$ws.Range("B7").Value = "PWD"
$ws.Range("C7").Value = 1

$ws.Range("A8").Value = "fra"
$ws.Range("B8").Value = "OTP"
$ws.Range("C8").Value = 2

$ws.Range("A9").Value = "fra"
$ws.Range("B9").Value = "FINGERPRINT"
$ws.Range("C9").Value = 3

$ws.Range("A10").Value = "fra"
$ws.Range("B10").Value = "IRIS"
$ws.Range("C10").Value = 4

$ws.Range("A11").Value = "fra"
$ws.Range("B11").Value = "FACE"
$ws.Range("C11").Value = 5

# D column ("is_active") must hold the text "TRUE" (shared string), not a
# boolean -- copy the existing "TRUE" text cell down instead of typing the
# literal, which Excel would otherwise auto-convert to a boolean.
for ($r = 7; $r -le 11; $r++) {
    $ws.Range("D2").Copy()
    $ws.Range("D$r").PasteSpecial(-4104)
}

# Match the final workbook's view state: active cell F9 and a slightly
# narrower default column width.
$ws.Range("F9").Select()
$ws.StandardWidth = 8.4296875
